$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 250 and 251 (existing rows 250-255 shift down to 252-257)
$ws.Range("A250:A251").EntireRow.Insert()

# New row 250 data
$ws.Range("A250").Value = 4
$ws.Range("B250").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C250").Value = "Los Lagos"
$ws.Range("D250").Value = 44448
$ws.Range("D250").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E250").Value = 10
$ws.Range("F250").Value = "Fruta"
$ws.Range("G250").Value = 100102
$ws.Range("H250").Value = "Cítricos"
$ws.Range("I250").Value = 100102003
$ws.Range("J250").Value = "Limón"
$ws.Range("K250").Value = "Sin especificar"
$ws.Range("L250").Value = "1a amarillo"
$ws.Range("M250").Value = 400
$ws.Range("N250").Value = 8500
$ws.Range("O250").Value = 8500
$ws.Range("P250").Value = 8500
$ws.Range("Q250").Value = "$/malla 16 kilos"
$ws.Range("R250").Value = "Región de O'Higgins"
$ws.Range("S250").Value = 531
$ws.Range("T250").Value = 16

# New row 251 data
$ws.Range("A251").Value = 4
$ws.Range("B251").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C251").Value = "Los Lagos"
$ws.Range("D251").Value = 44448
$ws.Range("D251").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E251").Value = 10
$ws.Range("F251").Value = "Fruta"
$ws.Range("G251").Value = 100102
$ws.Range("H251").Value = "Cítricos"
$ws.Range("I251").Value = 100102003
$ws.Range("J251").Value = "Limón"
$ws.Range("K251").Value = "Sin especificar"
$ws.Range("L251").Value = "2a amarillo"
$ws.Range("M251").Value = 200
$ws.Range("N251").Value = 7500
$ws.Range("O251").Value = 7500
$ws.Range("P251").Value = 7500
$ws.Range("Q251").Value = "$/malla 16 kilos"
$ws.Range("R251").Value = "Región de O'Higgins"
$ws.Range("S251").Value = 469
$ws.Range("T251").Value = 16
